$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protocol")

# Move the "Clients need timestamps..." comment to row 13 (F13), expanding the
# text with the extra sentence about measuring transport open time. This
# removes the old shared string (previously used nowhere else) and appends
# the new, longer string at the end of the shared string table.
$ws.Range("F13").Value = 'Clients need timestamps to determine if a streaming connection is being buffered, and can also use them as NOOPs to prevent a connection from closing. Or use them as an "please send me an initial message" to determine how long it really took to open the transport.'

# The longer comment needs a taller row to keep wrapping correctly.
$ws.Rows.Item(13).RowHeight = 25.5

# Reflect the author's last selection being on A13 after the edit.
$ws.Range("A13").Select() | Out-Null
